$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Job Number", "Customer Name", "Job State", "Job In Time", "Job Out Time", "Remark")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("D4").Select()
